$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF) - copy the formatting from the existing
# header cell H1 (bold, bordered, centered) so the new headers match the rest
# of the row, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-19
$data = @(
    @(3, 4),
    @(6, 7),
    @(7, 9),
    @(8, 9),
    @(7, 9),
    @(5, 6),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(4, 5),
    @(9, 9),
    @(6, 6),
    @(9, 9),
    @(3, 4),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
